$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record (email) is inserted as the new row 2, pushing the
# existing "jairist@aol.com" record down to row 3.
$ws.Rows.Item(2).Insert() | Out-Null

$ws.Cells.Item(2, 1).Value = "Login"
$ws.Cells.Item(2, 2).Value = "jairistasdfsdfdddasddssas@aol.com"
$ws.Cells.Item(2, 3).Value = "Monitorde15"

# Row insertion does not re-home the existing hyperlink, so drop the
# (now stale) hyperlink collection and recreate both links pointing at
# their correct cells.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 2), "mailto:jairist@aol.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "mailto:jairistasdfsdfdddasddssas@aol.com") | Out-Null

# Re-apply the Hyperlink cell style (Hyperlinks.Add alone creates a
# duplicate style entry instead of reusing the existing one).
$ws.Cells.Item(2, 2).Style = "Hyperlink"
$ws.Cells.Item(3, 2).Style = "Hyperlink"

# New "resultado" column header.
$ws.Cells.Item(1, 4).Value = "resultado"

# Match the saved selection state.
$ws.Range("C5").Select() | Out-Null
